# Version Final 03 fecha: 25/05/2023
#
# 1) Header row: replace underscores with spaces in a handful of column
#    titles (Razon_Social -> Razon Social, etc).
# 2) Body rows: blank out the placeholder "NA" text that used to sit in
#    several CURP / Primer_Apellido / Segundo_Apellido / Nombre_Comercial /
#    Fecha_Operacion / Estatus cells — these become empty cells instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header renames -----------------------------------------------
$ws.Range("C1").Value = "Razon Social"
$ws.Range("E1").Value = "Primer Apellido"
$ws.Range("F1").Value = "Segundo Apellido"
$ws.Range("G1").Value = "Nombre Comercial"
$ws.Range("H1").Value = "Fecha Operacion"

# --- 2) Clear stray "NA" placeholders ---------------------------------
$naCells = @(
    "D2", "E2", "F2",
    "D3", "E3", "F3",
    "D4", "E4", "F4",
    "D5", "E5", "F5", "G5",
    "D6", "E6", "F6", "H6", "I6",
    "D7", "E7", "F7",
    "D8", "E8", "F8",
    "G9",
    "G10",
    "G11",
    "G13",
    "G14"
)

foreach ($addr in $naCells) {
    $ws.Range($addr).Value = ""
}
